$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H32").Value = 1784.1578
$ws.Range("I32").Value = 2382.889
$ws.Range("J32").Value = 1245.3
$ws.Range("K32").Value = 2382.889
$ws.Range("L32").Value = 1245.3
$ws.Range("M32").Value = -2056.889
$ws.Range("N32").Value = -1897.3

$ws.Range("H33").Value = 359.6
$ws.Range("I33").Value = 361.90475
$ws.Range("J33").Value = 347.5
$ws.Range("K33").Value = 361.90475
$ws.Range("L33").Value = 347.5
$ws.Range("M33").Value = -132.90475
$ws.Range("N33").Value = -805.5

$ws.Range("H80").Value = 1462.75
$ws.Range("I80").Value = 1480
$ws.Range("J80").Value = 1434
$ws.Range("K80").Value = 4440
$ws.Range("L80").Value = 4302
$ws.Range("M80").Value = -3442
$ws.Range("N80").Value = -6298

$ws.Range("H83").Value = 1462.75
$ws.Range("I83").Value = 1480
$ws.Range("J83").Value = 1434
$ws.Range("K83").Value = 13320
$ws.Range("L83").Value = 12906
$ws.Range("M83").Value = -8328
$ws.Range("N83").Value = -22890

$ws.Range("H104").Value = 495.75
$ws.Range("I104").Value = 161
$ws.Range("J104").Value = 1500
$ws.Range("K104").Value = 483
$ws.Range("L104").Value = 4500
$ws.Range("M104").Value = 1264
$ws.Range("N104").Value = -7994

$ws.Range("H129").Value = 1006.1
$ws.Range("J129").Value = 1219.1428
$ws.Range("L129").Value = 3657.4284
$ws.Range("N129").Value = -13657.4284

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 8823.017
$ws.Range("I32").Value = 4863.854
$ws.Range("J32").Value = 23441.46
$ws.Range("K32").Value = 4863.854
$ws.Range("L32").Value = 23441.46
$ws.Range("M32").Value = -4576.854
$ws.Range("N32").Value = -24015.46

$ws.Range("H61").Value = 2590.4473
$ws.Range("I61").Value = 2057.2
$ws.Range("J61").Value = 4590.125
$ws.Range("K61").Value = 2057.2
$ws.Range("L61").Value = 4590.125
$ws.Range("M61").Value = -1845.2
$ws.Range("N61").Value = -5014.125

$ws.Range("H74").Value = 4437.724
$ws.Range("I74").Value = 533.0952
$ws.Range("J74").Value = 14687.375
$ws.Range("K74").Value = 533.0952
$ws.Range("L74").Value = 14687.375
$ws.Range("M74").Value = 340.9048
$ws.Range("N74").Value = -16435.375

$ws.Range("H77").Value = 4437.724
$ws.Range("I77").Value = 533.0952
$ws.Range("J77").Value = 14687.375
$ws.Range("K77").Value = 2665.476
$ws.Range("L77").Value = 73436.875
$ws.Range("M77").Value = 1702.524
$ws.Range("N77").Value = -82172.875

$ws.Range("H122").Value = 1911.2549
$ws.Range("I122").Value = 1653.5807
$ws.Range("J122").Value = 2310.65
$ws.Range("K122").Value = 4960.742099999999
$ws.Range("L122").Value = 6931.950000000001
$ws.Range("M122").Value = -2510.742099999999
$ws.Range("N122").Value = -11831.95

$ws.Range("H132").Value = 1997.3429
$ws.Range("I132").Value = 1790.2
$ws.Range("J132").Value = 2515.2
$ws.Range("K132").Value = 5370.6
$ws.Range("L132").Value = 7545.599999999999
$ws.Range("M132").Value = -2840.6
$ws.Range("N132").Value = -12605.6

$ws.Range("H136").Value = 2590.4473
$ws.Range("I136").Value = 2057.2
$ws.Range("J136").Value = 4590.125
$ws.Range("K136").Value = 6171.599999999999
$ws.Range("L136").Value = 13770.375
$ws.Range("M136").Value = -3621.599999999999
$ws.Range("N136").Value = -18870.375

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H29").Value = 1300
$ws.Range("I29").Value = 1300
$ws.Range("K29").Value = 1300
$ws.Range("M29").Value = -1011

$ws.Range("H134").Value = 32837.742
$ws.Range("I134").Value = 39069.207
$ws.Range("K134").Value = 117207.621
$ws.Range("M134").Value = -114672.621

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 1750.7778
$ws.Range("I31").Value = 1338.1
$ws.Range("J31").Value = 2266.625
$ws.Range("K31").Value = 1338.1
$ws.Range("L31").Value = 2266.625
$ws.Range("M31").Value = -1043.1
$ws.Range("N31").Value = -2856.625

$ws.Range("H34").Value = 1750.7778
$ws.Range("I34").Value = 1338.1
$ws.Range("J34").Value = 2266.625
$ws.Range("K34").Value = 1338.1
$ws.Range("L34").Value = 2266.625
$ws.Range("M34").Value = -1136.1
$ws.Range("N34").Value = -2670.625

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H113").Value = 439.0476
$ws.Range("I113").Value = 507.125
$ws.Range("J113").Value = 397.15384
$ws.Range("K113").Value = 1521.375
$ws.Range("L113").Value = 1191.46152
$ws.Range("M113").Value = 648.625
$ws.Range("N113").Value = -5531.46152

$ws.Range("H131").Value = 1787844.8
$ws.Range("J131").Value = 2042994.6
$ws.Range("L131").Value = 6128983.800000001
$ws.Range("N131").Value = -6139063.800000001

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H2").Value = 79.666664
$ws.Range("I2").Value = 41
$ws.Range("J2").Value = 99
$ws.Range("K2").Value = 41
$ws.Range("L2").Value = 99
$ws.Range("M2").Value = 72
$ws.Range("N2").Value = -325

$ws.Range("H5").Value = 2152.875
$ws.Range("J5").Value = 2190.0645
$ws.Range("L5").Value = 2190.0645
$ws.Range("N5").Value = -2414.0645

$ws.Range("H20").Value = 420833.34
$ws.Range("I20").Value = 457272.72
$ws.Range("J20").Value = 20000
$ws.Range("K20").Value = 457272.72
$ws.Range("L20").Value = 20000
$ws.Range("M20").Value = -457027.72
$ws.Range("N20").Value = -20490

$ws.Range("H44").Value = 9500
$ws.Range("I44").Value = 4000
$ws.Range("J44").Value = 15000
$ws.Range("K44").Value = 4000
$ws.Range("L44").Value = 15000
$ws.Range("M44").Value = -3404
$ws.Range("N44").Value = -16192

$ws.Range("H132").Value = 1844.8049
$ws.Range("I132").Value = 1772.1666
$ws.Range("J132").Value = 1901.6522
$ws.Range("K132").Value = 5316.4998
$ws.Range("L132").Value = 5704.9566
$ws.Range("M132").Value = -2786.4998
$ws.Range("N132").Value = -10764.9566

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 1285
$ws.Range("I7").Value = 1396.6666
$ws.Range("K7").Value = 1396.6666
$ws.Range("M7").Value = -1284.6666

$ws.Range("H55").Value = 285.7143
$ws.Range("I55").Value = 280.76923
$ws.Range("J55").Value = 350
$ws.Range("K55").Value = 280.76923
$ws.Range("L55").Value = 350
$ws.Range("M55").Value = -107.76923
$ws.Range("N55").Value = -696

$ws.Range("H61").Value = 2700.5
$ws.Range("I61").Value = 1220.8
$ws.Range("K61").Value = 1220.8
$ws.Range("M61").Value = -1018.8

$ws.Range("H113").Value = 2700.5
$ws.Range("I113").Value = 1220.8
$ws.Range("K113").Value = 1220.8
$ws.Range("M113").Value = 949.2

$ws.Range("H122").Value = 13277.556
$ws.Range("I122").Value = 14487.25
$ws.Range("K122").Value = 43461.75
$ws.Range("M122").Value = -41011.75

$ws.Range("H126").Value = 1285
$ws.Range("I126").Value = 1396.6666
$ws.Range("K126").Value = 4189.9998
$ws.Range("M126").Value = -1719.9998

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 1951.5714
$ws.Range("I122").Value = 2030.5
$ws.Range("J122").Value = 1920
$ws.Range("K122").Value = 6091.5
$ws.Range("L122").Value = 5760
$ws.Range("M122").Value = -3641.5
$ws.Range("N122").Value = -10660
